# Update "paises.xlsx" (sheet "Pais") with the latest COVID-19 snapshot.
# The underlying feed re-sorted a handful of neighbouring countries (by
# total cases) as their counts were refreshed, so some rows swap contents
# in addition to the plain numeric refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refresh timestamp caption (row 1, col A) ---------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 30 de Septiembre de 2020 a las 23:48"

# --- Plain numeric refresh (country / row order unchanged) ----------------
# row -> @{ col = value }
$plainUpdates = @{
    4   = @{ B = 7439355; C = 33209; D = 4683312; E = 2544467; G = 791;  H = 211576 }  # Estados Unidos
    13  = @{ B = 674339;  C = 1767;  D = 608112;  E = 49493;   G = 67;   H = 16734  }  # Sudafrica
    66  = @{ B = 46626;   C = 144;   D = 45757;   E = 568 }                            # Ghana
    85  = @{ B = 19724;   C = 55;    D = 19291;   E = 313 }                            # Costa de Marfil
    114 = @{ B = 7502;    C = 14;    D = 7120;    E = 221 }                            # Mauritania
    130 = @{ B = 4840;    C = 4;     D = 3154;    E = 1657 }                           # Ruanda
    137 = @{ B = 3963;    C = 29;    D = 3259;    E = 677;    G = 1;    H = 27    }    # Aruba
    147 = @{ B = 2894;    C = 48;    D = 1680;    E = 1134;   G = 2;    H = 80    }    # Guyana
    160 = @{ B = 1784;    C = 25;    D = 1348;    E = 388 }                            # Togo
}

foreach ($row in $plainUpdates.Keys) {
    $cols = $plainUpdates[$row]
    foreach ($col in $cols.Keys) {
        $colIndex = [int][char]$col - [int][char]'A' + 1
        $ws.Cells.Item($row, $colIndex).Value = $cols[$col]
    }
}

# --- Rank swaps: country data (incl. name) moves to a new row -------------
# Each entry is the FULL new row content (A..H) for the given row number.
$rowRewrites = @{
    153 = @("Burkina Faso",          2056, 24, 1335, 663, 0, 0,  58)
    154 = @("Principado de Andorra", 2050, 84, 1432, 565, 0, 0,  53)
    155 = @("Yemen",                 2034, 3,  1286, 161, 0, 0, 587)
    156 = @("Uruguay",               2033, 0,  1771, 214, 0, 0,  48)

    166 = @("Republica del Chad",    1200, 7,  1007, 108, 0, 0,  85)
    167 = @("Niger",                 1197, 1,  1114,  14, 0, 0,  69)

    215 = @("Montserrat",              13, 0,    12,   0, 0, 0,   1)
    216 = @("Islas Malvinas",          13, 0,    13,   0, 0, 0,   0)
}

foreach ($row in $rowRewrites.Keys) {
    $values = $rowRewrites[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}
